$wb = $excel.ActiveWorkbook

# Mapping for plot names and species names
$plotMap = @{ "样地1" = "A1"; "样地2" = "A2" }
$speciesMap = @{ "A" = "松树"; "B" = "杨树"; "C" = "柳树" }

# Sheet 1: 计算结果 - columns A (plot) and B (species), rows 2-6
$ws1 = $wb.Worksheets.Item("计算结果")
for ($r = 2; $r -le 6; $r++) {
    $aVal = $ws1.Cells.Item($r, 1).Value2
    if ($plotMap.ContainsKey($aVal)) {
        $ws1.Cells.Item($r, 1).Value2 = $plotMap[$aVal]
    }
    $bVal = $ws1.Cells.Item($r, 2).Value2
    if ($speciesMap.ContainsKey($bVal)) {
        $ws1.Cells.Item($r, 2).Value2 = $speciesMap[$bVal]
    }
}

# Sheet 2: 统计：物种 - column A (species), rows 2-4
$ws2 = $wb.Worksheets.Item("统计：物种")
for ($r = 2; $r -le 4; $r++) {
    $aVal = $ws2.Cells.Item($r, 1).Value2
    if ($speciesMap.ContainsKey($aVal)) {
        $ws2.Cells.Item($r, 1).Value2 = $speciesMap[$aVal]
    }
}

# Sheet 3: 统计：样地 - column A (plot), rows 2-3
$ws3 = $wb.Worksheets.Item("统计：样地")
for ($r = 2; $r -le 3; $r++) {
    $aVal = $ws3.Cells.Item($r, 1).Value2
    if ($plotMap.ContainsKey($aVal)) {
        $ws3.Cells.Item($r, 1).Value2 = $plotMap[$aVal]
    }
}

# Sheet 4: 统计：样地-物种 - columns A (plot) and B (species), rows 2-6
$ws4 = $wb.Worksheets.Item("统计：样地-物种")
for ($r = 2; $r -le 6; $r++) {
    $aVal = $ws4.Cells.Item($r, 1).Value2
    if ($plotMap.ContainsKey($aVal)) {
        $ws4.Cells.Item($r, 1).Value2 = $plotMap[$aVal]
    }
    $bVal = $ws4.Cells.Item($r, 2).Value2
    if ($speciesMap.ContainsKey($bVal)) {
        $ws4.Cells.Item($r, 2).Value2 = $speciesMap[$bVal]
    }
}
